$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 5; this shifts the former rows 5..47 down to 6..48
$ws.Rows("5:5").Insert()

# Fill the new row 5 ("0x000A" / "总花样数" / "16bit")
$ws.Range("A5").Value = "0x000A"
$ws.Range("B4").Value = "当前花样号"
$ws.Range("B5").Value = "总花样数"
$ws.Range("C5").Value = "16bit"

# Match the saved selection from the commit
$ws.Range("C16").Select()
